# Fix a typo in the "FindNewCarTest" worksheet ("Toyta Cars" -> "Toyota Cars")
# and update the worksheet's saved selection, matching the commit
# "added GetCarNameTest + some Readme comments".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindNewCarTest")
$ws.Activate()

# Correct the misspelled car name in C2 ("Toyta Cars" -> "Toyota Cars").
$ws.Range("C2").Value = "Toyota Cars"

# Update the active selection on the sheet to C9.
$ws.Range("C9").Select()
